$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 values: National/International, Status, Patent No
$ws.Range("B9").Value = "National"
$ws.Range("C9").Value = "Granted"
$ws.Range("D9").Value = 454732

# Widen column D to match column A's width (50.85546875 characters)
$ws.Columns.Item(4).ColumnWidth = 50

# Move the active selection to D2
$ws.Range("D2").Select()
